$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (Fernandez) - all rows below shift up by one
$ws.Rows.Item(2).Delete()

# Update recalculated values for rows 2-16 (post-shift) to match new period
# Row 2: Gomez
$ws.Range("A2").Value = "Gomez"
$ws.Range("B2").Value = 41456
$ws.Range("C2").Value = 12.09315068493151
$ws.Range("D2").Value = "Yoda-Sin Noches"
$ws.Range("E2").Value = 32
$ws.Range("F2").Value = 6.215593413784655
$ws.Range("G2").Value = 2.347299343765784
$ws.Range("H2").Value = 39
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2.666666666666667
$ws.Range("M2").Value = 2
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = 4.879012345679012
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0.3336076817558298

# Row 3: Bravo
$ws.Range("A3").Value = "Bravo"
$ws.Range("B3").Value = 41456
$ws.Range("C3").Value = 12.09315068493151
$ws.Range("D3").Value = "Yoda-Sin Noches"
$ws.Range("E3").Value = 31
$ws.Range("F3").Value = 6.215593413784655
$ws.Range("G3").Value = 2.179034157832754
$ws.Range("H3").Value = 35
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 2.666666666666667
$ws.Range("M3").Value = 2
$ws.Range("N3").Value = 2
$ws.Range("O3").Value = 4.360655737704918
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0.3322404371584699

# Row 4: Iñiguez
$ws.Range("A4").Value = "Iñiguez"
$ws.Range("B4").Value = 41640
$ws.Range("C4").Value = 11.58904109589041
$ws.Range("D4").Value = "Yoda-Sin Noches"
$ws.Range("E4").Value = 26
$ws.Range("F4").Value = 6.272659886921544
$ws.Range("G4").Value = 2.305233047282527
$ws.Range("H4").Value = 38
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 2.666666666666667
$ws.Range("M4").Value = 2
$ws.Range("N4").Value = 2
$ws.Range("O4").Value = 4.639357429718875
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0.3255689424364123

# Row 5: Breinbauer
$ws.Range("A5").Value = "Breinbauer"
$ws.Range("B5").Value = 41640
$ws.Range("C5").Value = 11.58904109589041
$ws.Range("D5").Value = "Yoda-Sin Noches"
$ws.Range("E5").Value = 28
$ws.Range("F5").Value = 6.272659886921544
$ws.Range("G5").Value = 2.305233047282527
$ws.Range("H5").Value = 38
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 2.666666666666667
$ws.Range("M5").Value = 2
$ws.Range("N5").Value = 2
$ws.Range("O5").Value = 4.676923076923077
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0.3282051282051282

# Row 6: Arredondo
$ws.Range("A6").Value = "Arredondo"
$ws.Range("B6").Value = 41852
$ws.Range("C6").Value = 11.00821917808219
$ws.Range("D6").Value = "Knight-Tardes"
$ws.Range("E6").Value = 29
$ws.Range("F6").Value = 6.338410388579262
$ws.Range("G6").Value = 7.916876998149114
$ws.Range("H6").Value = 35
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 16
$ws.Range("L6").Value = 4
$ws.Range("M6").Value = 2
$ws.Range("N6").Value = 2
$ws.Range("O6").Value = 4.32520325203252
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 1.977235772357723
$ws.Range("S6").Value = 0.4943089430894309

# Row 7: Carrasco
$ws.Range("A7").Value = "Carrasco"
$ws.Range("B7").Value = 41852
$ws.Range("C7").Value = 11.00821917808219
$ws.Range("D7").Value = "Knight-Tardes"
$ws.Range("E7").Value = 26
$ws.Range("F7").Value = 6.338410388579262
$ws.Range("G7").Value = 8.211341073531914
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 35
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 16
$ws.Range("L7").Value = 4
$ws.Range("M7").Value = 2
$ws.Range("N7").Value = 2
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 4.273092369477911
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 1.953413654618474
$ws.Range("S7").Value = 0.4883534136546185

# Row 8: Culaciati
$ws.Range("A8").Value = "Culaciati"
$ws.Range("B8").Value = 41852
$ws.Range("C8").Value = 11.00821917808219
$ws.Range("D8").Value = "Knight-Tardes"
$ws.Range("E8").Value = 29
$ws.Range("F8").Value = 6.338410388579262
$ws.Range("G8").Value = 8.211341073531914
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 35
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 16
$ws.Range("L8").Value = 4
$ws.Range("M8").Value = 2
$ws.Range("N8").Value = 2
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 4.32520325203252
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 1.977235772357723
$ws.Range("S8").Value = 0.4943089430894309

# Row 9: Contreras
$ws.Range("A9").Value = "Contreras"
$ws.Range("B9").Value = 42736
$ws.Range("C9").Value = 8.586301369863014
$ws.Range("D9").Value = "Knight-Tardes"
$ws.Range("E9").Value = 29
$ws.Range("F9").Value = 6.61257757473692
$ws.Range("G9").Value = 8.362779740871641
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 38
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 16
$ws.Range("L9").Value = 4
$ws.Range("M9").Value = 2
$ws.Range("N9").Value = 2
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 4.695934959349593
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 1.977235772357723
$ws.Range("S9").Value = 0.4943089430894309

# Row 10: Cisternas
$ws.Range("A10").Value = "Cisternas"
$ws.Range("B10").Value = 43040
$ws.Range("C10").Value = 7.753424657534246
$ws.Range("D10").Value = "Knight-Tardes"
$ws.Range("E10").Value = 14
$ws.Range("F10").Value = 6.706861312963085
$ws.Range("G10").Value = 8.41325929665155
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 39
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 16
$ws.Range("L10").Value = 4
$ws.Range("M10").Value = 2
$ws.Range("N10").Value = 2
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 4.542528735632184
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 1.86360153256705
$ws.Range("S10").Value = 0.4659003831417624

# Row 11: Pio
$ws.Range("A11").Value = "Pio"
$ws.Range("B11").Value = 43113
$ws.Range("C11").Value = 7.553424657534246
$ws.Range("D11").Value = "Padawan-Sin Fijo"
$ws.Range("E11").Value = 31
$ws.Range("F11").Value = 6.729501815892393
$ws.Range("G11").Value = 7.529867070503142
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 6
$ws.Range("K11").Value = 16
$ws.Range("L11").Value = 6.666666666666664
$ws.Range("M11").Value = 3
$ws.Range("N11").Value = 3
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0.7475409836065573
$ws.Range("R11").Value = 1.99344262295082
$ws.Range("S11").Value = 0.8306010928961745

# Row 12: Alvo
$ws.Range("A12").Value = "Alvo"
$ws.Range("B12").Value = 43770
$ws.Range("C12").Value = 5.753424657534246
$ws.Range("D12").Value = "Padawan-Sin Fijo"
$ws.Range("E12").Value = 28
$ws.Range("F12").Value = 6.933266342256174
$ws.Range("G12").Value = 7.866397442369204
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 6
$ws.Range("K12").Value = 17
$ws.Range("L12").Value = 6.666666666666664
$ws.Range("M12").Value = 3
$ws.Range("N12").Value = 3
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0.7384615384615385
$ws.Range("R12").Value = 2.092307692307692
$ws.Range("S12").Value = 0.8205128205128203

# Row 13: Boettiger
$ws.Range("A13").Value = "Boettiger"
$ws.Range("B13").Value = 44440
$ws.Range("C13").Value = 3.917808219178082
$ws.Range("D13").Value = "Padawan-Sin Fijo"
$ws.Range("E13").Value = 29
$ws.Range("F13").Value = 7.141062739004625
$ws.Range("G13").Value = 7.954736664984045
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 6
$ws.Range("K13").Value = 17
$ws.Range("L13").Value = 6.999999999999997
$ws.Range("M13").Value = 3
$ws.Range("N13").Value = 3
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0.7414634146341462
$ws.Range("R13").Value = 2.100813008130081
$ws.Range("S13").Value = 0.8650406504065037

# Row 14: Loch
$ws.Range("A14").Value = "Loch"
$ws.Range("B14").Value = 44713
$ws.Range("C14").Value = 3.16986301369863
$ws.Range("D14").Value = "Padawan-Sin Fijo"
$ws.Range("E14").Value = 28
$ws.Range("F14").Value = 7.225732017082725
$ws.Range("G14").Value = 8.556284704694631
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 6
$ws.Range("K14").Value = 18
$ws.Range("L14").Value = 7.999999999999996
$ws.Range("M14").Value = 3
$ws.Range("N14").Value = 4
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0.7384615384615385
$ws.Range("R14").Value = 2.215384615384616
$ws.Range("S14").Value = 0.9846153846153842

# Row 15: Rubio
$ws.Range("A15").Value = "Rubio"
$ws.Range("B15").Value = 44713
$ws.Range("C15").Value = 3.16986301369863
$ws.Range("D15").Value = "Padawan-Sin Fijo"
$ws.Range("E15").Value = 14
$ws.Range("F15").Value = 7.225732017082725
$ws.Range("G15").Value = 8.556284704694631
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = 18
$ws.Range("L15").Value = 7.999999999999996
$ws.Range("M15").Value = 4
$ws.Range("N15").Value = 3
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0.6988505747126437
$ws.Range("R15").Value = 2.096551724137931
$ws.Range("S15").Value = 0.9318007662835244

# Row 16: Recluta1
$ws.Range("A16").Value = "Recluta1"
$ws.Range("B16").Value = 45383
$ws.Range("C16").Value = 1.334246575342466
$ws.Range("D16").Value = "Padawan-Sin Fijo"
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = 7.433528413831178
$ws.Range("G16").Value = 9.284031633854989
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 8
$ws.Range("K16").Value = 19
$ws.Range("L16").Value = 8.999999999999998
$ws.Range("M16").Value = 4
$ws.Range("N16").Value = 4
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0.9353846153846154
$ws.Range("R16").Value = 2.221538461538461
$ws.Range("S16").Value = 1.052307692307692
